$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so that
# numeric-looking strings (e.g. "233.36", "1.00") are preserved exactly,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "93.407.71"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.415.29"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "233.36"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "621.26"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").Value = "1.45"
$ws.Range("E7").Value = "  +5.73%  "
$ws.Range("D8").Value = "0.392"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "0.988"
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("D11").Value = "3.416.99"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "43.80"
$ws.Range("E12").Value = "  +8.30%  "
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("E14").Value = "  +4.73%  "
$ws.Range("D15").Value = "93.229.35"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "4.052.99"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "8.37"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("D19").Value = "3.416.51"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "18.22"
$ws.Range("E20").Value = "  +7.94%  "
$ws.Range("D21").Value = "11.68"
$ws.Range("E21").Value = "  +7.19%  "
$ws.Range("E22").Value = "  +12.36%  "
$ws.Range("D23").Value = "3.38"
$ws.Range("E23").Value = "  +7.83%  "
$ws.Range("D24").Value = "499.15"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "6.75"
$ws.Range("E25").Value = "  +9.09%  "
$ws.Range("D26").Value = "0.0000182"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "86.62"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").Value = "12.05"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.595.03"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "11.40"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.141"
$ws.Range("E31").Value = "  +6.27%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "2.72"
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "0.177"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "0.554"
$ws.Range("E36").Value = "  +5.25%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "29.11"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "555.53"
$ws.Range("E38").Value = "  +6.19%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "7.49"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.41"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.149"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.899"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "23.71"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "1.70"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0419"
$ws.Range("E46").Value = "  +6.13%  "
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").Value = "3.62"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "5.50"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "8.15"
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "53.13"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "2.10"
$ws.Range("E51").Value = "  -2.60%  "
